# Reorganización completa: limpieza de módulos antiguos, nuevas entregas y optimización
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "condicion_corporal"

# Clear old formatting (bold white font on blue fill) from the previous header row
$ws.Cells.ClearFormats()

# Reset the custom column widths that existed on columns A-F
for ($c = 1; $c -le 6; $c++) {
    $ws.Columns.Item($c).ColumnWidth = 8.43
}

# New header row: codigo, descripcion, puntuacion, escala, especie,
# caracteristicas, recomendaciones, estado
$ws.Range("A1").Value = "codigo"
$ws.Range("B1").Value = "descripcion"
$ws.Range("C1").Value = "puntuacion"
$ws.Range("D1").Value = "escala"
$ws.Range("E1").Value = "especie"
$ws.Range("F1").Value = "caracteristicas"
$ws.Range("G1").Value = "recomendaciones"
$ws.Range("H1").Value = "estado"

Write-Host "done"
